# Malawi specs workbook: "Added ability to incorporate urban centers and clusters"
#
# Updates the rural/urban electrification-ratio inputs (columns AM/AN on
# row 2 - the lower-case "rural_elec_ratio" / "urban_elec_ratio" fields)
# to the recalculated values that account for urban centers/clusters, and
# leaves the cursor/selection where the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# rural_elec_ratio
$ws.Range("AM2").Value = 0.46741966083849035
# urban_elec_ratio
$ws.Range("AN2").Value = 0.56986553015578867

# Scroll the view over towards the right-hand columns and leave the
# selection on AN7, matching where the sheet was left after the edit.
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 22
$win.ScrollRow = 1
$ws.Range("AN7").Select() | Out-Null
